# ---------------------------------------------------------------------------
# busi448 / 10_leverage.xlsx -- add a "Repo" sheet (shorting / repo-leverage
# notebook) next to the existing "Margin" sheet, with its own assumptions
# block, two balance-sheet tables, and named ranges used by the formulas.
# ---------------------------------------------------------------------------

$xlPasteFormats = -4122
$xlCenter = -4108
$xlRight  = -4152

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original sheet, add the new one right after it.
# ---------------------------------------------------------------------------
$margin = $wb.Worksheets.Item(1)
$margin.Name = "Margin"

$repo = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $margin)
$repo.Name = "Repo"

# ---------------------------------------------------------------------------
# 2. Column widths for Repo (A narrow index col, B/D labels, C/E numbers).
# ---------------------------------------------------------------------------
$repo.Columns.Item(1).ColumnWidth = 4.36328125
$repo.Columns.Item(2).ColumnWidth = 16
$repo.Columns.Item(3).ColumnWidth = 10.90625
$repo.Columns.Item(4).ColumnWidth = 16
$repo.Columns.Item(5).ColumnWidth = 10.90625

# ---------------------------------------------------------------------------
# 3. Section headers -- copy the bold "section title" look from Margin!A1.
# ---------------------------------------------------------------------------
$margin.Range("A1").Copy() | Out-Null
$repo.Range("A1").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("A1").Value = "Assumptions"

$margin.Range("A1").Copy() | Out-Null
$repo.Range("A7").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("A7").Value = "Initial Balance Sheet (buy bond and borrow in repo market)"

$margin.Range("A1").Copy() | Out-Null
$repo.Range("A16").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("A16").Value = "Balance Sheet after Realized Return (prior to returning cash)"

$margin.Range("A1").Copy() | Out-Null
$repo.Range("A17").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------------
# 4. Assumptions block (rows 2-5).
# ---------------------------------------------------------------------------
$repo.Range("B2").Value = "MV"
$repo.Range("D2").Value = "Implied haircut"
$repo.Range("B3").Value = "Initial cash"
$repo.Range("D3").Value = "Implied repo rate"
$repo.Range("B4").Value = "Repurchase price"
$repo.Range("B5").Value = "Term (days)"

# Bold/blue "input" cells (C2:C5) -- copy look from Margin!C10 (bold blue)
# then drop its percent number format back to General.
$margin.Range("C10").Copy() | Out-Null
$repo.Range("C2:C5").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("C2:C5").NumberFormat = "General"

$repo.Range("C2").Value = 1000
$repo.Range("C3").Value = 980
$repo.Range("C4").Value = 981
$repo.Range("C5").Value = 10

# Implied haircut / repo rate -- same look as Margin!C7 (0.0%, no font override).
$margin.Range("C7").Copy() | Out-Null
$repo.Range("E2:E3").PasteSpecial($xlPasteFormats) | Out-Null

$repo.Range("E2").Formula = "=1-init_cash/MV"
$repo.Range("E3").Formula = "=(repo_price/init_cash-1)*360/term"

$repo.Range("F2").Formula = "=FORMULATEXT(haircut)"
$repo.Range("F3").Formula = "=FORMULATEXT(repo_rate)"

# ---------------------------------------------------------------------------
# 5. First balance sheet (rows 8-11) -- copy the whole Assets/Liab table
#    layout from Margin!B2:E5 (keeps borders/number formats identical).
# ---------------------------------------------------------------------------
$margin.Range("B2:E5").Copy() | Out-Null
$repo.Range("B8").PasteSpecial($xlPasteFormats) | Out-Null

$repo.Range("B8").Value = "Assets"
$repo.Range("D8").Value = "Liabilities & Equity"

$repo.Range("B9").Value = "Bond (repo'd)"
$repo.Range("C9").Formula = "=MV"
$repo.Range("D9").Value = "Repo (Cash Loan)"
$repo.Range("E9").Formula = "=init_cash"

$repo.Range("D10").Value = "Equity"
$repo.Range("E10").Formula = "=C11-E9"

$repo.Range("C11").Formula = "=SUM(C9:C10)"
$repo.Range("E11").Formula = "=SUM(E9:E10)"

$repo.Rows.Item(11).RowHeight = 15

# ---------------------------------------------------------------------------
# 6. Leverage / percent margin rows (12-15).
# ---------------------------------------------------------------------------
$repo.Range("B12").HorizontalAlignment = $xlCenter
$repo.Range("B15").HorizontalAlignment = $xlCenter

$margin.Range("C3").Copy() | Out-Null
$repo.Range("C12").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("C15").PasteSpecial($xlPasteFormats) | Out-Null

$margin.Range("E3").Copy() | Out-Null
$repo.Range("E12").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("E12").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
$repo.Range("E12:E15").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("E13").Copy() | Out-Null
$repo.Range("E13,E14,E15").PasteSpecial($xlPasteFormats) | Out-Null

$repo.Range("B13").Value = "Leverage"
$repo.Range("B13").HorizontalAlignment = $xlRight
$repo.Range("C13").Formula = "=E9/E10"
$repo.Range("C13").NumberFormat = "0%"
$repo.Range("D13").Formula = "=FORMULATEXT(C13)"

$repo.Range("B14").Value = "Percent margin"
$repo.Range("B14").HorizontalAlignment = $xlRight
$margin.Range("C7").Copy() | Out-Null
$repo.Range("C14").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("C14").Formula = "=E10/C11"
$repo.Range("D14").Formula = "=FORMULATEXT(C14)"

# ---------------------------------------------------------------------------
# 7. "Bond return (total)" assumption + second balance sheet (rows 17-21).
# ---------------------------------------------------------------------------
$repo.Range("B17").Value = "Bond return (total)"

$margin.Range("C10").Copy() | Out-Null
$repo.Range("C17").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("C17").NumberFormat = "0.0%"
$repo.Range("C17").Value = 0.005
$repo.Rows.Item(17).RowHeight = 15

$margin.Range("B2:E5").Copy() | Out-Null
$repo.Range("B18").PasteSpecial($xlPasteFormats) | Out-Null

$repo.Range("B18").Value = "Assets"
$repo.Range("D18").Value = "Liabilities & Equity"

$repo.Range("B19").Value = "Bond (repo'd)"
$repo.Range("C19").Formula = "=C9*(1+return)"
$repo.Range("D19").Value = "Repo (Cash Loan)"
$repo.Range("E19").Formula = "=repo_price"

$repo.Range("D20").Value = "Equity"
$repo.Range("E20").Formula = "=C21-E19"

$repo.Range("C21").Formula = "=SUM(C19:C20)"
$repo.Range("E21").Formula = "=SUM(E19:E20)"

$repo.Rows.Item(21).RowHeight = 15

# ---------------------------------------------------------------------------
# 8. Leverage / percent margin / levered return rows (22-27).
# ---------------------------------------------------------------------------
$repo.Range("B22").HorizontalAlignment = $xlCenter
$margin.Range("C3").Copy() | Out-Null
$repo.Range("C22").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("E13").Copy() | Out-Null
$repo.Range("E22,E23,E24,E25,E26,E27").PasteSpecial($xlPasteFormats) | Out-Null

$repo.Range("B23").Value = "Leverage"
$repo.Range("B23").HorizontalAlignment = $xlRight
$repo.Range("C23").Formula = "=E19/E20"
$repo.Range("C23").NumberFormat = "0%"
$repo.Range("D23").Formula = "=FORMULATEXT(C23)"

$repo.Range("B24").Value = "Percent margin"
$repo.Range("B24").HorizontalAlignment = $xlRight
$margin.Range("C7").Copy() | Out-Null
$repo.Range("C24").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("C24").Formula = "=E20/C21"
$repo.Range("D24").Formula = "=FORMULATEXT(C24)"

$repo.Range("B25").Value = "Levered return"
$repo.Range("B25").HorizontalAlignment = $xlRight
$margin.Range("C7").Copy() | Out-Null
$repo.Range("C25").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("C25").Formula = "=E20/E10-1"
$repo.Range("D25").Formula = "=FORMULATEXT(C25)"

$repo.Range("B26").Value = "Levered return (formula)"
$repo.Range("B26").HorizontalAlignment = $xlRight
$repo.Range("B26").WrapText = $true
$margin.Range("C7").Copy() | Out-Null
$repo.Range("C26").PasteSpecial($xlPasteFormats) | Out-Null
$repo.Range("C26").WrapText = $true
$repo.Range("C26").Formula = "=return+leverage*(return-repo_rate*term/360)"
$repo.Range("D26").Formula = "=FORMULATEXT(C26)"
$repo.Rows.Item(26).RowHeight = 29

$repo.Range("B27").HorizontalAlignment = $xlRight
$repo.Range("C13").Copy() | Out-Null
$repo.Range("C27,D27").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------------
# 9. Named ranges used throughout the formulas above.
# ---------------------------------------------------------------------------
$wb.Names.Add("MV", $repo.Range("C2"))
$wb.Names.Add("init_cash", $repo.Range("C3"))
$wb.Names.Add("repo_price", $repo.Range("C4"))
$wb.Names.Add("term", $repo.Range("C5"))
$wb.Names.Add("haircut", $repo.Range("E2"))
$wb.Names.Add("repo_rate", $repo.Range("E3"))
$wb.Names.Add("leverage", $repo.Range("C13"))
$wb.Names.Add("return", $repo.Range("C17"))

# ---------------------------------------------------------------------------
# 10. View state: zoom + re-select A1 on both sheets.
# ---------------------------------------------------------------------------
$repo.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 205

$margin.Activate() | Out-Null
$margin.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 280

Write-Output "done"
